$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status text for existing rows 2 and 3 (row 4 keeps its original status)
$ws.Range("L2").Value = "מחשבים שלא טופלו"
$ws.Range("L3").Value = "מחשבים שטופלו ונלקחו"

# Append a new record in row 5
$ws.Range("A5").Value = "אלדד עזוז"

# Column B holds a date-looking value that must stay plain text (matches
# the existing rows, which store it as inline text, not a real date).
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2023-10-17"
$ws.Range("B5").Style = "Normal"

$ws.Range("C5").Value = "הכנסה לדומיין"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "2UA24222"
$ws.Range("F5").Value = 1212
$ws.Range("G5").Value = "בסמח"
$ws.Range("H5").Value = 9996333
$ws.Range("I5").Value = "eldad@gmail.com"
$ws.Range("J5").Value = 502361254
$ws.Range("K5").Value = "אין"
$ws.Range("L5").Value = "מחשבים שטופלו ונלקחו"
$ws.Range("M5").Value = "ACYDBNiG4l0vuauSyQwWZACPpzh8iDZLIfRlzvKrqVv4HgeEZEIz5MhWbVMoxJVy6XDDDis"
